$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 33, shifting existing rows 33:84 down to 34:85
$ws.Rows.Item(33).Insert()

# Populate the new row 33 with the latest September transaction entry
$ws.Range("R33").Value = "balance your axis"
$ws.Range("S33").Value = "2024-09-06 09:55:31"
